# Auto-generated edit script applying scheduled-runner updates to Hades_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 15569.357
$ws.Range("I88").Value = 1125.75
$ws.Range("J88").Value = 21346.8
$ws.Range("K88").Value = 1125.75
$ws.Range("L88").Value = 21346.8
$ws.Range("M88").Value = -719.75
$ws.Range("N88").Value = -22158.8

# Row 91
$ws.Range("H91").Value = 15569.357
$ws.Range("I91").Value = 1125.75
$ws.Range("J91").Value = 21346.8
$ws.Range("K91").Value = 1125.75
$ws.Range("L91").Value = 21346.8
$ws.Range("M91").Value = 278.25
$ws.Range("N91").Value = -24154.8

# Row 96
$ws.Range("H96").Value = 1049.6666
$ws.Range("I96").Value = 791.5
$ws.Range("J96").Value = 1143.5454
$ws.Range("K96").Value = 2374.5
$ws.Range("L96").Value = 3430.6362
$ws.Range("M96").Value = -1001.5
$ws.Range("N96").Value = -6176.6362

# Row 132
$ws.Range("H132").Value = 4459887
$ws.Range("I132").Value = 5875.6
$ws.Range("J132").Value = 49000000
$ws.Range("K132").Value = 17626.8
$ws.Range("L132").Value = 147000000
$ws.Range("M132").Value = -15096.8
$ws.Range("N132").Value = -147005060

# Row 138
$ws.Range("H138").Value = 2566113.8
$ws.Range("I138").Value = 1368.5883
$ws.Range("J138").Value = 3474461
$ws.Range("K138").Value = 4105.7649
$ws.Range("L138").Value = 10423383
$ws.Range("M138").Value = 1034.2351
$ws.Range("N138").Value = -10433663

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4175.7188
$ws.Range("I45").Value = 4206.9614
$ws.Range("J45").Value = 4040.3333
$ws.Range("K45").Value = 4206.9614
$ws.Range("L45").Value = 4040.3333
$ws.Range("M45").Value = -3829.9614
$ws.Range("N45").Value = -4794.3333

# Row 61
$ws.Range("H61").Value = 52738532
$ws.Range("I61").Value = 66735150
$ws.Range("J61").Value = 251200
$ws.Range("K61").Value = 66735150
$ws.Range("L61").Value = 251200
$ws.Range("M61").Value = -66734938
$ws.Range("N61").Value = -251624

# Row 80
$ws.Range("H80").Value = 44147
$ws.Range("J80").Value = 47196
$ws.Range("L80").Value = 47196
$ws.Range("N80").Value = -49192

# Row 83
$ws.Range("H83").Value = 44147
$ws.Range("J83").Value = 47196
$ws.Range("L83").Value = 141588
$ws.Range("N83").Value = -151572

# Row 122
$ws.Range("H122").Value = 1047.7878
$ws.Range("I122").Value = 1040.742
$ws.Range("K122").Value = 3122.226
$ws.Range("M122").Value = -672.2259999999997

# Row 136
$ws.Range("H136").Value = 52738532
$ws.Range("I136").Value = 66735150
$ws.Range("J136").Value = 251200
$ws.Range("K136").Value = 200205450
$ws.Range("L136").Value = 753600
$ws.Range("M136").Value = -200202900
$ws.Range("N136").Value = -758700

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Row 139
$ws.Range("H139").Value = 59881.5
$ws.Range("J139").Value = 59881.5
$ws.Range("L139").Value = 59881.5
$ws.Range("N139").Value = -70161.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1198.7142
$ws.Range("I94").Value = 1122.4546
$ws.Range("J94").Value = 1478.3334
$ws.Range("K94").Value = 1122.4546
$ws.Range("L94").Value = 1478.3334
$ws.Range("M94").Value = -671.4546
$ws.Range("N94").Value = -2380.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3101.16
$ws.Range("I31").Value = 2017
$ws.Range("J31").Value = 3195.4348
$ws.Range("K31").Value = 2017
$ws.Range("L31").Value = 3195.4348
$ws.Range("M31").Value = -1722
$ws.Range("N31").Value = -3785.4348

# Row 34
$ws.Range("H34").Value = 3101.16
$ws.Range("I34").Value = 2017
$ws.Range("J34").Value = 3195.4348
$ws.Range("K34").Value = 2017
$ws.Range("L34").Value = 3195.4348
$ws.Range("M34").Value = -1815
$ws.Range("N34").Value = -3599.4348

# Row 124
$ws.Range("H124").Value = 32000
$ws.Range("J124").Value = 32000
$ws.Range("L124").Value = 32000
$ws.Range("N124").Value = -36910

# Row 132
$ws.Range("H132").Value = 50973.24
$ws.Range("I132").Value = 2889.2856
$ws.Range("J132").Value = 147141.14
$ws.Range("K132").Value = 8667.856800000001
$ws.Range("L132").Value = 441423.42
$ws.Range("M132").Value = -6137.856800000001
$ws.Range("N132").Value = -446483.42

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Range("H94").Value = 3305.2632
$ws.Range("J94").Value = 3305.2632
$ws.Range("L94").Value = 9915.7896
$ws.Range("N94").Value = -11267.7896

# Row 113
$ws.Range("H113").Value = 601.1667
$ws.Range("I113").Value = 502.09677
$ws.Range("J113").Value = 688.9143
$ws.Range("K113").Value = 1506.29031
$ws.Range("L113").Value = 2066.7429
$ws.Range("M113").Value = 663.7096900000001
$ws.Range("N113").Value = -6406.7429

# Row 131
$ws.Range("H131").Value = 916.2817
$ws.Range("I131").Value = 504.14285
$ws.Range("J131").Value = 961.3594000000001
$ws.Range("K131").Value = 1512.42855
$ws.Range("L131").Value = 2884.0782
$ws.Range("M131").Value = 3527.57145
$ws.Range("N131").Value = -12964.0782

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 26324642
$ws.Range("I80").Value = 13498
$ws.Range("J80").Value = 55559244
$ws.Range("K80").Value = 13498
$ws.Range("L80").Value = 55559244
$ws.Range("M80").Value = -12500
$ws.Range("N80").Value = -55561240

# Row 83
$ws.Range("H83").Value = 26324642
$ws.Range("I83").Value = 13498
$ws.Range("J83").Value = 55559244
$ws.Range("K83").Value = 67490
$ws.Range("L83").Value = 277796220
$ws.Range("M83").Value = -62498
$ws.Range("N83").Value = -277806204

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2026
$ws.Range("I93").Value = 1900
$ws.Range("J93").Value = 2152
$ws.Range("K93").Value = 1900
$ws.Range("L93").Value = 2152
$ws.Range("M93").Value = -652
$ws.Range("N93").Value = -4648

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5707.2856
$ws.Range("I81").Value = 6975.5
$ws.Range("J81").Value = 5200
$ws.Range("K81").Value = 13951
$ws.Range("L81").Value = 10400
$ws.Range("M81").Value = -12890
$ws.Range("N81").Value = -12522

# Row 84
$ws.Range("H84").Value = 5707.2856
$ws.Range("I84").Value = 6975.5
$ws.Range("J84").Value = 5200
$ws.Range("K84").Value = 69755
$ws.Range("L84").Value = 52000
$ws.Range("M84").Value = -64451
$ws.Range("N84").Value = -62608
